$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "18:11 03-Dec-23"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "22166101"
$ws.Range("D2").Value = "thầy Long đỉnh vaixi luon "
